$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.708015918731689
$ws.Range("B1").Value = 2.616949796676636
$ws.Range("C1").Value = 2.263574600219727
$ws.Range("D1").Value = 1.49393904209137
$ws.Range("E1").Value = 0.935169517993927
